# Cronograma-July-Tesis.xlsx — "Actualizado a Fecha 23 de marzo Actualizado"
#
# Updates the progress (%) and notes column for the second block of use
# cases (rows 18-29), adding a "revisar nuevamente" note to each, and
# moves the current selection/viewport down to the newly-updated rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$note = "revisar nuevamente"

# Row -> new progress value (column D). Column E (merged E:I) gets the note.
$updates = @{
    18 = 0.8
    19 = 0.4
    20 = 0.4
    21 = 0.3
    22 = 0.5
    23 = 0.4
    24 = 0.5
    25 = 0.5
    26 = 0.5
    27 = 0.5
    28 = 0
    29 = 0.5
}

foreach ($row in $updates.Keys | Sort-Object) {
    $ws.Cells.Item($row, 4).Value = $updates[$row]
    $ws.Cells.Item($row, 5).Value = $note
}

# Recalculate dependent formulas (e.g. the D30 average) before saving.
$excel.Calculate()

# Move the visible selection to reflect where work continued.
$ws.Range("E27:I27").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
